$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '19.951.26'
$ws.Range("E2").Value = '  -5.09%  '
$ws.Range("D3").Value = '1.415.60'
$ws.Range("E3").Value = '  -5.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.69%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '276.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3666'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3094'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.032'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06530'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.491'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.190'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("D16").Value = '1.414.82'
$ws.Range("E16").Value = '  -6.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001018'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.37%  '
$ws.Range("E18").Value = '  -13.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -12.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.614'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.241'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.21%  '
$ws.Range("D25").Value = '19.961.74'
$ws.Range("E25").Value = '  -5.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.265'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '132.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.24%  '
$ws.Range("D29").Value = '1.571.17'
$ws.Range("E29").Value = '  -6.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.885'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -18.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.254'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8159'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -13.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07689'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.486'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.301'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.920'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05777'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9971'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02049'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.45'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1884'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.094'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.10%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.23%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5307'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.537'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5182'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '115.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.767'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.032'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.68%  '
$ws.Range("E51").Value = '  -0.61%  '
